$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting of the already-styled (but empty) row 59
#     down onto row 60, which currently carries a different / older
#     style set, before we populate both rows with data. ---
$ws.Range("A59:K59").Copy()
$ws.Range("A60:K60").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 59: ESRI entry ---
$ws.Range("A59").Value = "GIS Mapping Software, Location Intelligence & Spatial Analytics"
$ws.Range("B59").Value = "ESRI"
$ws.Range("C59").Value = 2024
$ws.Range("D59").Value = "ESRI"
$ws.Range("E59").Value = "Software/Package"
$ws.Range("F59").Value = "NA"
$ws.Range("G59").Value = "ESRI base map citation"
$ws.Range("H59").Value = "GIS"
$ws.Range("I59").Value = "yes"
$ws.Range("J59").Value = "yes"
$ws.Range("K59").Value = "yes"

# --- Row 60: QGIS entry ---
$ws.Range("A60").Value = "QGIS Geographic Information System [Computer software]"
$ws.Range("B60").Value = "QGIS"
$ws.Range("C60").Value = 2024
$ws.Range("D60").Value = "QGIS developpment team"
$ws.Range("E60").Value = "Software/Package"
$ws.Range("F60").Value = "NA"
$ws.Range("G60").Value = "QGIS software citation"
$ws.Range("H60").Value = "GIS"
$ws.Range("I60").Value = "yes"
$ws.Range("J60").Value = "yes"
$ws.Range("K60").Value = "yes"

# --- Recalculate summary COUNTIF formulas that depend on the new rows ---
$wb.Application.Calculate()

# --- Update the view: scroll position + active selection, mirroring the
#     author's navigation to the newly-added rows. ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D62").Select()
